# Auto-applied edit corresponding to the Betfair Back/Lay games-of-the-day update
# for 2025-11-17: the Algerian Ligue 1 fixture (ES Mostaganem vs USM Alger) is
# removed, every remaining fixture shifts up one row, and all odds are refreshed
# with newly scraped values; a new Argentinian Primera Division fixture
# (CA Platense vs Gimnasia La Plata) slots into the now-vacant final row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete first data row (row 2); everything below shifts up by one
# row and the sheet dimension/used range shrinks from A1:AO9 to A1:AO8 automatically.
$ws.Rows.Item(2).Delete() | Out-Null

# Keep the Date/Time columns as plain text instead of letting Excel auto-parse
# the assigned strings into date/time serial numbers.
$ws.Range("B2:C8").NumberFormat = "@"

$numRows = 7
$numCols = 41
$data = New-Object 'object[,]' $numRows,$numCols
# Row 2 (source row 2)
$data[0,0] = "Italian Serie C"
$data[0,1] = "2025-11-17"
$data[0,2] = "16:30:00"
$data[0,3] = "Foggia"
$data[0,4] = "Cavese 1919"
$data[0,5] = 3.5
$data[0,6] = 3.95
$data[0,7] = 2.3
$data[0,8] = 2.44
$data[0,9] = 3.1
$data[0,10] = 3.35
$data[0,11] = 1.57
$data[0,12] = 1.11
$data[0,13] = 2.82
$data[0,14] = 1.53
$data[0,15] = 1.6
$data[0,16] = 2.54
$data[0,17] = 1.21
$data[0,18] = 5.1
$data[0,19] = 2.04
$data[0,20] = 1.81
$data[0,21] = 1.69
$data[0,22] = 1.35
$data[0,23] = 10.5
$data[0,24] = 7.6
$data[0,25] = 17
$data[0,26] = 65
$data[0,27] = 10.5
$data[0,28] = 7.2
$data[0,29] = 12.5
$data[0,30] = 55
$data[0,31] = 26
$data[0,32] = 17.5
$data[0,33] = 23
$data[0,34] = 90
$data[0,35] = 190
$data[0,36] = 150
$data[0,37] = 120
$data[0,38] = 1000
$data[0,39] = 160
$data[0,40] = 60
# Row 3 (source row 3)
$data[1,0] = "Italian Serie C"
$data[1,1] = "2025-11-17"
$data[1,2] = "16:30:00"
$data[1,3] = "Pergolettese"
$data[1,4] = "Giana Erminio"
$data[1,5] = 2.66
$data[1,6] = 2.88
$data[1,7] = 3.05
$data[1,8] = 3.35
$data[1,9] = 2.96
$data[1,10] = 3.15
$data[1,11] = 1.63
$data[1,12] = 1.14
$data[1,13] = 2.56
$data[1,14] = 1.59
$data[1,15] = 1.48
$data[1,16] = 2.86
$data[1,17] = 1.16
$data[1,18] = 6.2
$data[1,19] = 2.2
$data[1,20] = 1.72
$data[1,21] = 1.43
$data[1,22] = 1.53
$data[1,23] = 8.199999999999999
$data[1,24] = 8.800000000000001
$data[1,25] = 26
$data[1,26] = 100
$data[1,27] = 7.8
$data[1,28] = 7.2
$data[1,29] = 18.5
$data[1,30] = 85
$data[1,31] = 17
$data[1,32] = 21
$data[1,33] = 38
$data[1,34] = 150
$data[1,35] = 130
$data[1,36] = 120
$data[1,37] = 120
$data[1,38] = 1000
$data[1,39] = 1000
$data[1,40] = 600
# Row 4 (source row 4)
$data[2,0] = "Spanish Segunda Division"
$data[2,1] = "2025-11-17"
$data[2,2] = "16:30:00"
$data[2,3] = "Leonesa"
$data[2,4] = "Malaga"
$data[2,5] = 2.36
$data[2,6] = 2.42
$data[2,7] = 3.6
$data[2,8] = 3.75
$data[2,9] = 3.2
$data[2,10] = 3.25
$data[2,11] = 1.52
$data[2,12] = 1.11
$data[2,13] = 3.05
$data[2,14] = 1.46
$data[2,15] = 1.7
$data[2,16] = 2.4
$data[2,17] = 1.25
$data[2,18] = 4.6
$data[2,19] = 1.99
$data[2,20] = 1.97
$data[2,21] = 1.36
$data[2,22] = 1.71
$data[2,23] = 9.800000000000001
$data[2,24] = 11
$data[2,25] = 25
$data[2,26] = 80
$data[2,27] = 8
$data[2,28] = 6.8
$data[2,29] = 15
$data[2,30] = 55
$data[2,31] = 13.5
$data[2,32] = 11.5
$data[2,33] = 20
$data[2,34] = 70
$data[2,35] = 36
$data[2,36] = 29
$data[2,37] = 48
$data[2,38] = 140
$data[2,39] = 28
$data[2,40] = 70
# Row 5 (source row 5)
$data[3,0] = "Argentinian Primera Division"
$data[3,1] = "2025-11-17"
$data[3,2] = "17:00:00"
$data[3,3] = "Barracas Central"
$data[3,4] = "Huracan"
$data[3,5] = 3.15
$data[3,6] = 3.3
$data[3,7] = 3.05
$data[3,8] = 3.1
$data[3,9] = 2.74
$data[3,10] = 2.78
$data[3,11] = 1.81
$data[3,12] = 1.21
$data[3,13] = 2.12
$data[3,14] = 1.84
$data[3,15] = 1.35
$data[3,16] = 3.65
$data[3,17] = 1.11
$data[3,18] = 9
$data[3,19] = 2.52
$data[3,20] = 1.55
$data[3,21] = 1.47
$data[3,22] = 1.44
$data[3,23] = 6
$data[3,24] = 7.2
$data[3,25] = 18
$data[3,26] = 120
$data[3,27] = 7.4
$data[3,28] = 6.8
$data[3,29] = 17
$data[3,30] = 65
$data[3,31] = 18.5
$data[3,32] = 17.5
$data[3,33] = 34
$data[3,34] = 120
$data[3,35] = 70
$data[3,36] = 70
$data[3,37] = 130
$data[3,38] = 300
$data[3,39] = 110
$data[3,40] = 1000
# Row 6 (source row 6)
$data[4,0] = "Argentinian Primera Division"
$data[4,1] = "2025-11-17"
$data[4,2] = "17:00:00"
$data[4,3] = "Belgrano"
$data[4,4] = "Union Santa Fe"
$data[4,5] = 2.24
$data[4,6] = 2.3
$data[4,7] = 4
$data[4,8] = 4.1
$data[4,9] = 3.1
$data[4,10] = 3.25
$data[4,11] = 1.63
$data[4,12] = 1.14
$data[4,13] = 2.58
$data[4,14] = 1.6
$data[4,15] = 1.51
$data[4,16] = 2.86
$data[4,17] = 1.17
$data[4,18] = 6.2
$data[4,19] = 2.22
$data[4,20] = 1.74
$data[4,21] = 1.32
$data[4,22] = 1.77
$data[4,23] = 8.4
$data[4,24] = 11
$data[4,25] = 27
$data[4,26] = 110
$data[4,27] = 6.6
$data[4,28] = 7.2
$data[4,29] = 18.5
$data[4,30] = 75
$data[4,31] = 12.5
$data[4,32] = 11.5
$data[4,33] = 27
$data[4,34] = 110
$data[4,35] = 32
$data[4,36] = 34
$data[4,37] = 70
$data[4,38] = 230
$data[4,39] = 40
$data[4,40] = 390
# Row 7 (source row 7)
$data[5,0] = "Argentinian Primera Division"
$data[5,1] = "2025-11-17"
$data[5,2] = "17:00:00"
$data[5,3] = "Defensa y Justicia"
$data[5,4] = "Independiente Rivadavia"
$data[5,5] = 1.94
$data[5,6] = 1.96
$data[5,7] = 5
$data[5,8] = 5.3
$data[5,9] = 3.35
$data[5,10] = 3.45
$data[5,11] = 1.53
$data[5,12] = 1.11
$data[5,13] = 3.05
$data[5,14] = 1.46
$data[5,15] = 1.68
$data[5,16] = 2.4
$data[5,17] = 1.24
$data[5,18] = 4.7
$data[5,19] = 2.06
$data[5,20] = 1.84
$data[5,21] = 1.23
$data[5,22] = 2.04
$data[5,23] = 10.5
$data[5,24] = 14.5
$data[5,25] = 95
$data[5,26] = 140
$data[5,27] = 7.2
$data[5,28] = 7.6
$data[5,29] = 21
$data[5,30] = 80
$data[5,31] = 10.5
$data[5,32] = 10.5
$data[5,33] = 34
$data[5,34] = 100
$data[5,35] = 21
$data[5,36] = 24
$data[5,37] = 48
$data[5,38] = 180
$data[5,39] = 20
$data[5,40] = 1000
# Row 8 (source row 8)
$data[6,0] = "Argentinian Primera Division"
$data[6,1] = "2025-11-17"
$data[6,2] = "19:30:00"
$data[6,3] = "CA Platense"
$data[6,4] = "Gimnasia La Plata"
$data[6,5] = 2.66
$data[6,6] = 2.7
$data[6,7] = 3.3
$data[6,8] = 3.45
$data[6,9] = 3
$data[6,10] = 3.1
$data[6,11] = 1.69
$data[6,12] = 1.15
$data[6,13] = 2.36
$data[6,14] = 1.71
$data[6,15] = 1.44
$data[6,16] = 3.05
$data[6,17] = 1.14
$data[6,18] = 7
$data[6,19] = 2.34
$data[6,20] = 1.67
$data[6,21] = 1.4
$data[6,22] = 1.59
$data[6,23] = 7.6
$data[6,24] = 8.6
$data[6,25] = 21
$data[6,26] = 80
$data[6,27] = 7.2
$data[6,28] = 7.2
$data[6,29] = 16.5
$data[6,30] = 65
$data[6,31] = 15
$data[6,32] = 13
$data[6,33] = 29
$data[6,34] = 260
$data[6,35] = 44
$data[6,36] = 44
$data[6,37] = 85
$data[6,38] = 250
$data[6,39] = 60
$data[6,40] = 390

$ws.Range("A2:AO8").Value2 = $data

Write-Host "Applied update: used range now" $ws.UsedRange.Address()
